$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "C2"
$ws.Range("D2").Value = "N3"
$ws.Range("E2").Value = "C4"
$ws.Range("F2").Value = "C6"
$ws.Range("G2").Value = "C12"
$ws.Range("H2").Value = "N13"
$ws.Range("I2").Value = "C5"
$ws.Range("J2").Value = "C4"
$ws.Range("C3").Value = "C2"
$ws.Range("D3").Value = "N3"
$ws.Range("E3").Value = "C4"
$ws.Range("F3").Value = "C6"
$ws.Range("G3").Value = "C12"
$ws.Range("H3").Value = "N13"
$ws.Range("I3").Value = "C5"
$ws.Range("J3").Value = "C4"
$ws.Range("C4").Value = "C2"
$ws.Range("D4").Value = "N3"
$ws.Range("E4").Value = "C4"
$ws.Range("F4").Value = "C6"
$ws.Range("G4").Value = "C12"
$ws.Range("H4").Value = "N13"
$ws.Range("I4").Value = "C5"
$ws.Range("J4").Value = "C4"
$ws.Range("C5").Value = "C2"
$ws.Range("D5").Value = "N3"
$ws.Range("E5").Value = "C4"
$ws.Range("F5").Value = "C6"
$ws.Range("G5").Value = "C12"
$ws.Range("H5").Value = "N13"
$ws.Range("I5").Value = "C5"
$ws.Range("J5").Value = "C4"
$ws.Range("C6").Value = "C7"
$ws.Range("D6").Value = "N6"
$ws.Range("E6").Value = "C4"
$ws.Range("F6").Value = "C3"
$ws.Range("G6").Value = "C2"
$ws.Range("H6").Value = "N15"
$ws.Range("I6").Value = "C8"
$ws.Range("J6").Value = "C7"
$ws.Range("C7").Value = "C7"
$ws.Range("D7").Value = "N15"
$ws.Range("E7").Value = "C2"
$ws.Range("F7").Value = "C3"
$ws.Range("G7").Value = "C4"
$ws.Range("H7").Value = "N6"
$ws.Range("I7").Value = "C8"
$ws.Range("J7").Value = "C7"
$ws.Range("C8").Value = "C2"
$ws.Range("D8").Value = "N3"
$ws.Range("E8").Value = "C4"
$ws.Range("F8").Value = "C5"
$ws.Range("G8").Value = "C8"
$ws.Range("H8").Value = "N9"
$ws.Range("I8").Value = "C6"
$ws.Range("J8").Value = "C5"
$ws.Range("C9").Value = "C2"
$ws.Range("D9").Value = "N9"
$ws.Range("E9").Value = "C8"
$ws.Range("F9").Value = "C5"
$ws.Range("G9").Value = "C4"
$ws.Range("H9").Value = "N3"
$ws.Range("I9").Value = "C6"
$ws.Range("J9").Value = "C5"
$ws.Range("C10").Value = "C4"
$ws.Range("D10").Value = "N3"
$ws.Range("E10").Value = "C2"
$ws.Range("F10").Value = "C10"
$ws.Range("G10").Value = "C8"
$ws.Range("H10").Value = "N7"
$ws.Range("I10").Value = "C9"
$ws.Range("J10").Value = "C8"
$ws.Range("C11").Value = "C9"
$ws.Range("D11").Value = "N5"
$ws.Range("E11").Value = "C4"
$ws.Range("F11").Value = "C3"
$ws.Range("G11").Value = "C6"
$ws.Range("H11").Value = "N8"
$ws.Range("I11").Value = "C7"
$ws.Range("J11").Value = "C6"
$ws.Range("C12").Value = "C9"
$ws.Range("D12").Value = "N5"
$ws.Range("E12").Value = "C4"
$ws.Range("F12").Value = "C3"
$ws.Range("G12").Value = "C6"
$ws.Range("H12").Value = "N8"
$ws.Range("I12").Value = "C7"
$ws.Range("J12").Value = "C6"
$ws.Range("C13").Value = "C5"
$ws.Range("D13").Value = "N4"
$ws.Range("E13").Value = "C3"
$ws.Range("F13").Value = "C2"
$ws.Range("G13").Value = "C8"
$ws.Range("H13").Value = "N7"
$ws.Range("I13").Value = "C1"
$ws.Range("J13").Value = "C2"
$ws.Range("C14").Value = "C2"
$ws.Range("D14").Value = "N3"
$ws.Range("E14").Value = "C4"
$ws.Range("F14").Value = "C5"
$ws.Range("G14").Value = "C7"
$ws.Range("H14").Value = "N10"
$ws.Range("I14").Value = "C6"
$ws.Range("J14").Value = "C5"
